$d = $word.ActiveDocument

# Rename the "IB Tabela 1" custom table style to "StatsTLF Tabela 1",
# matching the already-renamed "StatsTLF ..." family of styles in this
# template (StatsTLF Titulo 1/2, StatsTLF Normal 1, etc.).
$tableStyle = $d.Styles("IB Tabela 1")
$tableStyle.NameLocal = "StatsTLF Tabela 1"
